# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview sheet (zh-cn/de-de status columns) and on each language sheet's
#   Status column.
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   to reflect the new handoff generation.
# - Narrow columns E/F on Overview and column C on the zh-cn/de-de sheets (status
#   columns), which no longer need to be as wide for the shorter status text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-21 03:02:41"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-21 03:02:37"
$wsZh.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-21 03:02:41"
$wsDe.Columns.Item(3).ColumnWidth = 16.33
